$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the next day's log entry as a new row (row 99), matching the
# existing sheet's layout: Date, day-of-week, hour, ranking.
$row = 99

# Column A holds a date formatted as literal text ("2025/10/13"), matching
# every prior row in this log. Force text entry (NumberFormat "@") so Excel
# doesn't auto-convert it to a date serial, then restore the default
# "Normal" cell style so no extra formatting is left behind on the cell.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025/10/13"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = "月"
$ws.Cells.Item($row, 3).Value = 16
$ws.Cells.Item($row, 4).Value = 143
